$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 556
$ws1.Range("F6").Value = 512
$ws1.Range("F10").Value = 6818
$ws1.Range("F11").Value = 237
$ws1.Range("F13").Value = 3139
$ws1.Range("F14").Value = 207
$ws1.Range("F15").Value = 366
$ws1.Range("F17").Value = 554
$ws1.Range("F18").Value = 14

# Sheet "全部类型" (All types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value = 556
$ws4.Range("F8").Value = 512
$ws4.Range("F13").Value = 6818
$ws4.Range("F15").Value = 237
$ws4.Range("F17").Value = 3139
$ws4.Range("F18").Value = 207
$ws4.Range("F19").Value = 366
$ws4.Range("F21").Value = 554
$ws4.Range("F22").Value = 14
